# Apply updated crypto price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.175.18"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "3.431.93"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D5").Value = "'407.03"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").Value = "'133.04"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("D7").Value = "'0.591"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.674"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").Value = "'0.122"
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("D11").Value = "'42.53"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "3.938.18"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "'8.44"
$ws.Range("E14").Value = "  -3.48%  "
$ws.Range("D15").Value = "'19.90"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").Value = "3.408.98"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").Value = "62.092.04"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("D19").Value = "'10.99"
$ws.Range("D20").Value = "'0.0000131"
$ws.Range("E20").Value = "  -5.07%  "
$ws.Range("D21").Value = "'3.22"
$ws.Range("E21").Value = "  -4.72%  "
$ws.Range("D22").Value = "'84.64"
$ws.Range("E22").Value = "  +2.90%  "
$ws.Range("D23").Value = "'317.32"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "'12.82"
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").Value = "'3.13"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "'4.78"
$ws.Range("E26").Value = "  +8.94%  "
$ws.Range("D27").Value = "'29.79"
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("D28").Value = "'8.26"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "'7.73"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").Value = "'2.71"
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("D31").Value = "'0.173"
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'42.49"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'11.42"
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "'0.0483"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "'0.997"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "'3.44"
$ws.Range("E39").Value = "  -4.03%  "
$ws.Range("D40").Value = "'3.01"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").Value = "'138.75"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").Value = "'2.00"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").Value = "'3.98"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "'16.77"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").Value = "'21.36"
$ws.Range("E48").Value = "  -3.58%  "
$ws.Range("D49").Value = "2.132.27"
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("E50").Value = "  -5.58%  "
$ws.Range("D51").Value = "'1.88"
$ws.Range("E51").Value = "  +1.65%  "
